{"js": "// Replace each two-digit multiplication problem's text with its new value.\n// Mapping derived from the authoritative diff: every \"<w:t>\" run that holds\n// a problem like \"19\u00d741=\" is swapped for a new problem, e.g. \"54\u00d711=\".\nconst replacements = [\n  [\"19\u00d741=\", \"54\u00d711=\"],\n  [\"15\u00d756=\", \"37\u00d737=\"],\n  [\"18\u00d720=\", \"15\u00d794=\"],\n  [\"53\u00d732=\", \"16\u00d731=\"],\n  [\"43\u00d776=\", \"23\u00d765=\"],\n  [\"21\u00d789=\", \"94\u00d791=\"],\n  [\"73\u00d749=\", \"58\u00d780=\"],\n  [\"92\u00d732=\", \"89\u00d727=\"],\n  [\"89\u00d777=\", \"32\u00d730=\"],\n  [\"68\u00d714=\", \"82\u00d714=\"],\n  [\"13\u00d768=\", \"43\u00d754=\"],\n  [\"42\u00d778=\", \"70\u00d743=\"],\n  [\"25\u00d771=\", \"50\u00d780=\"],\n  [\"38\u00d773=\", \"97\u00d736=\"],\n  [\"71\u00d761=\", \"61\u00d787=\"],\n  [\"74\u00d768=\", \"37\u00d740=\"],\n  [\"39\u00d775=\", \"54\u00d763=\"],\n  [\"63\u00d740=\", \"34\u00d725=\"],\n  [\"17\u00d736=\", \"35\u00d744=\"],\n  [\"97\u00d756=\", \"78\u00d769=\"],\n  [\"39\u00d742=\", \"44\u00d760=\"],\n  [\"82\u00d763=\", \"77\u00d743=\"],\n  [\"74\u00d778=\", \"98\u00d766=\"],\n  [\"49\u00d718=\", \"36\u00d796=\"],\n  [\"75\u00d758=\", \"85\u00d774=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication problem's text with its new value.\n# Mapping derived from the authoritative diff: every Find.Text value that\n# holds a problem like \"19\u00d741=\" is swapped for a new problem, e.g. \"54\u00d711=\".\n$replacements = @(\n    @{ Old = \"19\u00d741=\"; New = \"54\u00d711=\" },\n    @{ Old = \"15\u00d756=\"; New = \"37\u00d737=\" },\n    @{ Old = \"18\u00d720=\"; New = \"15\u00d794=\" },\n    @{ Old = \"53\u00d732=\"; New = \"16\u00d731=\" },\n    @{ Old = \"43\u00d776=\"; New = \"23\u00d765=\" },\n    @{ Old = \"21\u00d789=\"; New = \"94\u00d791=\" },\n    @{ Old = \"73\u00d749=\"; New = \"58\u00d780=\" },\n    @{ Old = \"92\u00d732=\"; New = \"89\u00d727=\" },\n    @{ Old = \"89\u00d777=\"; New = \"32\u00d730=\" },\n    @{ Old = \"68\u00d714=\"; New = \"82\u00d714=\" },\n    @{ Old = \"13\u00d768=\"; New = \"43\u00d754=\" },\n    @{ Old = \"42\u00d778=\"; New = \"70\u00d743=\" },\n    @{ Old = \"25\u00d771=\"; New = \"50\u00d780=\" },\n    @{ Old = \"38\u00d773=\"; New = \"97\u00d736=\" },\n    @{ Old = \"71\u00d761=\"; New = \"61\u00d787=\" },\n    @{ Old = \"74\u00d768=\"; New = \"37\u00d740=\" },\n    @{ Old = \"39\u00d775=\"; New = \"54\u00d763=\" },\n    @{ Old = \"63\u00d740=\"; New = \"34\u00d725=\" },\n    @{ Old = \"17\u00d736=\"; New = \"35\u00d744=\" },\n    @{ Old = \"97\u00d756=\"; New = \"78\u00d769=\" },\n    @{ Old = \"39\u00d742=\"; New = \"44\u00d760=\" },\n    @{ Old = \"82\u00d763=\"; New = \"77\u00d743=\" },\n    @{ Old = \"74\u00d778=\"; New = \"98\u00d766=\" },\n    @{ Old = \"49\u00d718=\"; New = \"36\u00d796=\" },\n    @{ Old = \"75\u00d758=\"; New = \"85\u00d774=\" }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
